# TimeSpentUC.xlsx - "final presentation + MVC + some rename"
#
# Adds a small "Calculated time" vs. "Real time spent" comparison table
# (rows 52-56) below the existing Function-Point-Calculation sheet, reusing
# the four 2nd-semester UC labels (rows 17-20) and rounding their computed
# "Calculated time (hours)" (column G) values to one decimal place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Base formatting for the new block -------------------------------------
# Style of I9 (plain cell with border, no fill) -> reused for the whole
# A52:C56 block (matches cellXfs index 2 in the authored file).
$ws.Range("I9").Copy()
$ws.Range("A52:C56").PasteSpecial(-4122)

# Style of A17 (bordered + coloured header fill) -> reused for the UC-name
# cells A53:A56, matching the look of A17:A20 above (cellXfs index 8).
$ws.Range("A17").Copy()
$ws.Range("A53:A56").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Header row --------------------------------------------------------
$ws.Range("B52").Value = "Calculated time"
$ws.Range("C52").Value = "Real time spent"

# --- Data rows: reuse the same UC labels as rows 17-20 ------------------
$ws.Range("A53").Value = $ws.Range("A17").Value()
$ws.Range("B53").Value = 38.9

$ws.Range("A54").Value = $ws.Range("A18").Value()
$ws.Range("B54").Value = 8

$ws.Range("A55").Value = $ws.Range("A19").Value()
$ws.Range("B55").Value = 11.3

$ws.Range("A56").Value = $ws.Range("A20").Value()
$ws.Range("B56").Value = 5

# --- View state: scroll down to the new table and select F55:F56 --------
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
$ws.Range("F55:F56").Select()

# --- Minor column-width tweak on column C (user nudged it slightly) -----
$ws.Columns.Item(3).ColumnWidth = 9.5
